$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Insert two new bulleted paragraphs ("Create Spring application from
#    spring.io" / "Follow Backend part. ...") right before the existing
#    "Front End:" paragraph (paragraph 2), using a numbered list (numId 6,
#    a brand-new list definition).
# ---------------------------------------------------------------------

# 1a. Mint a brand-new numbered-list definition by applying the default
#     numbered list to a throwaway paragraph appended at the very end of
#     the document, then delete that paragraph again. This leaves a new
#     <w:num>/<w:abstractNum> pair behind in numbering.xml that later
#     paragraphs can reference by numId.
$endRange = $d.Content
$endRange.Collapse(0)
$endRange.InsertParagraphAfter()
$mint = $d.Paragraphs($d.Paragraphs.Count)
$mint.Range.ListFormat.ApplyNumberDefault()
$mint2 = $d.Paragraphs($d.Paragraphs.Count)
$mint2.Range.Delete() | Out-Null

$frontEnd = $d.Paragraphs(2)
$insRange = $d.Range($frontEnd.Range.Start, $frontEnd.Range.End)

$newParasXml = @"
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
<pkg:xmlData>
<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
<w:body>
<w:p>
<w:pPr>
<w:pStyle w:val="ListParagraph"/>
<w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr>
<w:rPr><w:b/><w:bCs/><w:highlight w:val="yellow"/></w:rPr>
</w:pPr>
<w:r><w:t xml:space="preserve">Create </w:t></w:r>
<w:r><w:t xml:space="preserve">Spring application from </w:t></w:r>
<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>spring.io</w:t></w:r>
</w:p>
<w:p>
<w:pPr>
<w:pStyle w:val="ListParagraph"/>
<w:numPr><w:ilvl w:val="0"/><w:numId w:val="6"/></w:numPr>
<w:rPr><w:b/><w:bCs/><w:highlight w:val="yellow"/></w:rPr>
</w:pPr>
<w:r><w:t xml:space="preserve">Follow </w:t></w:r>
<w:r><w:rPr><w:b/><w:bCs/></w:rPr><w:t>Backend</w:t></w:r>
<w:r><w:t xml:space="preserve"> part. (Front End section is only for reference)</w:t></w:r>
</w:p>
<w:p w14:paraId="4DD0BB50" w14:textId="654A3B78" w:rsidR="00146E3D" w:rsidRPr="007761A0" w:rsidRDefault="003430F7">
<w:pPr><w:rPr><w:b/><w:bCs/></w:rPr></w:pPr>
<w:r w:rsidRPr="003430F7"><w:rPr><w:b/><w:bCs/><w:highlight w:val="yellow"/></w:rPr><w:t xml:space="preserve">Front </w:t></w:r>
<w:r w:rsidR="007761A0" w:rsidRPr="003430F7"><w:rPr><w:b/><w:bCs/><w:highlight w:val="yellow"/></w:rPr><w:t>End:</w:t></w:r>
</w:p>
</w:body>
</w:document>
</pkg:xmlData>
</pkg:part>
</pkg:package>
"@

$insRange.InsertXML($newParasXml)

Write-Output "Step 1 done. Paragraph count: $($d.Paragraphs.Count)"
